$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = -0.058712873528375502
$ws.Range("B1").Value = 0.058712872106111209
$ws.Range("A2").Value = 0.0050900061247426973
$ws.Range("B2").Value = -0.0050900075934149549
$ws.Range("A3").Value = 0.0081973339691865446
$ws.Range("B3").Value = -0.0081973355188688091

$ws.Columns.Item(1).ColumnWidth = 14.7109375
$ws.Columns.Item(2).ColumnWidth = 15.42578125
